$d = $word.ActiveDocument

# Grab the last paragraph in the document (the "Changes in local repo on
# 2304 docx" paragraph) and add a brand-new paragraph right after it,
# mirroring that paragraph's formatting (including its tab stop at 4221),
# then fill it in with the new sentence.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$newPara.Range.Text = "Second change on 2304 docx"
